$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.26599999999999
$ws.Range("B9").Value = 6.569899999999998
$ws.Range("C9").Value = -11.35130000000001
$ws.Range("D9").Value = -7.279399999999999
$ws.Range("C11").Value = -13.57160000000001
$ws.Range("B18").Value = 6.253299999999997
$ws.Range("B20").Value = 8.979799999999999
$ws.Range("C23").Value = -12.2416
$ws.Range("C24").Value = -13.0495
$ws.Range("C26").Value = -12.7993
$ws.Range("B27").Value = 6.436900000000009
$ws.Range("D27").Value = -7.773699999999999
$ws.Range("D29").Value = -7.1842
$ws.Range("D32").Value = -7.05309999999999
$ws.Range("C34").Value = -12.30340000000001
$ws.Range("B35").Value = 8.480200000000002
$ws.Range("C35").Value = -13.39620000000001
$ws.Range("D37").Value = -7.5297
$ws.Range("D38").Value = -7.242300000000001
$ws.Range("D41").Value = -8.939799999999993
$ws.Range("D45").Value = -7.155499999999996
$ws.Range("C48").Value = -10.73349999999999
$ws.Range("C49").Value = -13.94029999999999
$ws.Range("D51").Value = -8.561800000000003
$ws.Range("C52").Value = -11.1466
$ws.Range("D57").Value = -8.493599999999999
$ws.Range("D64").Value = -7.29859999999999
$ws.Range("C66").Value = -10.7851
$ws.Range("C67").Value = -10.6071
$ws.Range("B69").Value = 5.462599999999992
$ws.Range("B76").Value = 5.529199999999999
$ws.Range("B78").Value = 10.24110000000001
$ws.Range("C78").Value = -13.6074
$ws.Range("C80").Value = -13.26440000000001
$ws.Range("B82").Value = 5.1068
$ws.Range("D82").Value = -8.349000000000009
$ws.Range("B83").Value = 5.5549
$ws.Range("B93").Value = 4.863899999999996
$ws.Range("D93").Value = -7.173299999999994
$ws.Range("C99").Value = -13.2422
$ws.Range("D102").Value = -7.503299999999997
$ws.Range("C104").Value = -12.81090000000001
$ws.Range("D105").Value = -7.706300000000002
